$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (losing significant trailing zeros versus the literal target string),
# so force them to Text format before assignment.
foreach ($ref in @("D12", "D25", "D46", "D47")) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.884.47"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.895.13"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "0.7828"
$ws.Range("D6").Value = "243.72"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.3135"
$ws.Range("D9").Value = "25.77"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("D10").Value = "0.07353"
$ws.Range("E10").Value = "  +5.34%  "
$ws.Range("D11").Value = "0.08091"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "0.7730"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "5.505"
$ws.Range("E13").Value = "  +4.85%  "
$ws.Range("D14").Value = "1.903.34"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "93.92"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "6.214"
$ws.Range("E16").Value = "  +5.70%  "
$ws.Range("D17").Value = "29.869.70"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "13.98"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("D19").Value = "247.14"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("D20").Value = "0.000007823"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "8.155"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "2.129.36"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "0.1590"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").Value = "9.441"
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("D27").Value = "163.86"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").Value = "18.73"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "2.024"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").Value = "1.436"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("D31").Value = "1.543"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "4.487"
$ws.Range("E32").Value = "  +2.81%  "
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "4.062"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "1.239"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").Value = "0.7523"
$ws.Range("E36").Value = "  +3.04%  "
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "0.01938"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").Value = "1.140.97"
$ws.Range("E41").Value = "  +12.40%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.4474"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "74.36"
$ws.Range("E43").Value = "  +3.44%  "
$ws.Range("D44").Value = "5.959"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("D45").Value = "0.8526"
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("D46").Value = "1.0000"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "1.890"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").Value = "102.55"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").Value = "3.056"
$ws.Range("E49").Value = "  +6.17%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "7.524"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.771"
$ws.Range("E51").Value = "  -0.84%  "
